$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 26 (idx 0)
$ws.Range("H26").Value = 6000
$ws.Range("J26").Value = 6000
$ws.Range("L26").Value = 6000
$ws.Range("N26").Value = -6688

# Row 40 (idx 1)
$ws.Range("H40").Value = 899.5
$ws.Range("J40").Value = 899.5
$ws.Range("L40").Value = 899.5
$ws.Range("N40").Value = -1249.5

# Row 80 (idx 2)
$ws.Range("H80").Value = 1468.8
$ws.Range("I80").Value = 1468.8
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4406.4
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3408.4

# Row 83 (idx 3)
$ws.Range("H83").Value = 1468.8
$ws.Range("I83").Value = 1468.8
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13219.2
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -8227.199999999999

# Row 125 (idx 4)
$ws.Range("H125").Value = 9496
$ws.Range("J125").Value = 8999
$ws.Range("L125").Value = 80991
$ws.Range("N125").Value = -85911

$ws = $wb.Worksheets.Item("ARM")
# Row 24 (idx 5)
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0

# Row 61 (idx 6)
$ws.Range("H61").Value = 6749.6665
$ws.Range("I61").Value = 7199.6
$ws.Range("K61").Value = 7199.6
$ws.Range("M61").Value = -6987.6

# Row 74 (idx 7)
$ws.Range("H74").Value = 3666.5454
$ws.Range("I74").Value = 933
$ws.Range("J74").Value = 5228.5713
$ws.Range("K74").Value = 933
$ws.Range("L74").Value = 5228.5713
$ws.Range("M74").Value = -59
$ws.Range("N74").Value = -6976.5713

# Row 77 (idx 8)
$ws.Range("H77").Value = 3666.5454
$ws.Range("I77").Value = 933
$ws.Range("J77").Value = 5228.5713
$ws.Range("K77").Value = 4665
$ws.Range("L77").Value = 26142.8565
$ws.Range("M77").Value = -297
$ws.Range("N77").Value = -34878.85649999999

# Row 100 (idx 9)
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

# Row 136 (idx 10)
$ws.Range("H136").Value = 6749.6665
$ws.Range("I136").Value = 7199.6
$ws.Range("K136").Value = 21598.8
$ws.Range("M136").Value = -19048.8

$ws = $wb.Worksheets.Item("BSM")
# Row 106 (idx 11)
$ws.Range("H106").Value = 5869.25
$ws.Range("J106").Value = 5869.25
$ws.Range("L106").Value = 5869.25
$ws.Range("N106").Value = -8393.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (idx 12)
$ws.Range("H16").Value = 6145.143
$ws.Range("I16").Value = 6203.2
$ws.Range("J16").Value = 6000
$ws.Range("K16").Value = 6203.2
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = -5916.2
$ws.Range("N16").Value = -6574

# Row 31 (idx 13)
$ws.Range("H31").Value = 4279.647
$ws.Range("I31").Value = 1289.3077
$ws.Range("J31").Value = 13998.25
$ws.Range("K31").Value = 1289.3077
$ws.Range("L31").Value = 13998.25
$ws.Range("M31").Value = -994.3077000000001
$ws.Range("N31").Value = -14588.25

# Row 34 (idx 14)
$ws.Range("H34").Value = 4279.647
$ws.Range("I34").Value = 1289.3077
$ws.Range("J34").Value = 13998.25
$ws.Range("K34").Value = 1289.3077
$ws.Range("L34").Value = 13998.25
$ws.Range("M34").Value = -1087.3077
$ws.Range("N34").Value = -14402.25

# Row 86 (idx 15)
$ws.Range("H86").Value = 29783.727
$ws.Range("I86").Value = 8159
$ws.Range("J86").Value = 55733.4
$ws.Range("K86").Value = 8159
$ws.Range("L86").Value = 55733.4
$ws.Range("M86").Value = -7036
$ws.Range("N86").Value = -57979.4

# Row 89 (idx 16)
$ws.Range("H89").Value = 29783.727
$ws.Range("I89").Value = 8159
$ws.Range("J89").Value = 55733.4
$ws.Range("K89").Value = 40795
$ws.Range("L89").Value = 278667
$ws.Range("M89").Value = -35179
$ws.Range("N89").Value = -289899

# Row 113 (idx 17)
$ws.Range("H113").Value = 6145.143
$ws.Range("I113").Value = 6203.2
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 6203.2
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -4033.2
$ws.Range("N113").Value = -10340

# Row 132 (idx 18)
$ws.Range("H132").Value = 3748.4
$ws.Range("I132").Value = 2174.75
$ws.Range("J132").Value = 4797.5
$ws.Range("K132").Value = 6524.25
$ws.Range("L132").Value = 14392.5
$ws.Range("M132").Value = -3994.25
$ws.Range("N132").Value = -19452.5

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (idx 19)
$ws.Range("H4").Value = 125097416
$ws.Range("I4").Value = 75028890
$ws.Range("K4").Value = 225086670
$ws.Range("M4").Value = -225086558

# Row 34 (idx 20)
$ws.Range("H34").Value = 3395.4285
$ws.Range("I34").Value = 1900
$ws.Range("J34").Value = 3993.6
$ws.Range("K34").Value = 5700
$ws.Range("L34").Value = 11980.8
$ws.Range("M34").Value = -5616
$ws.Range("N34").Value = -12148.8

# Row 139 (idx 21)
$ws.Range("H139").Value = 5000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = 15000
$ws.Range("N139").Value = -25280

$ws = $wb.Worksheets.Item("GSM")
# Row 18 (idx 22)
$ws.Range("H18").Value = 15634.667
$ws.Range("I18").Value = 5000
$ws.Range("K18").Value = 5000
$ws.Range("M18").Value = -4707

# Row 23 (idx 23)
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("N23").Value = 0

# Row 132 (idx 24)
$ws.Range("H132").Value = 4270.3
$ws.Range("I132").Value = 3386.4285
$ws.Range("J132").Value = 6332.6665
$ws.Range("K132").Value = 10159.2855
$ws.Range("L132").Value = 18997.9995
$ws.Range("M132").Value = -7629.2855
$ws.Range("N132").Value = -24057.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (idx 25)
$ws.Range("H16").Value = 1148.8889
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# Row 38 (idx 26)
$ws.Range("H38").Value = 4900
$ws.Range("J38").Value = 4900
$ws.Range("L38").Value = 4900
$ws.Range("N38").Value = -5720

# Row 82 (idx 27)
$ws.Range("H82").Value = 1500
$ws.Range("J82").Value = 1500
$ws.Range("L82").Value = 1500
$ws.Range("N82").Value = -2222

# Row 85 (idx 28)
$ws.Range("H85").Value = 1500
$ws.Range("J85").Value = 1500
$ws.Range("L85").Value = 1500
$ws.Range("N85").Value = -3996

$ws = $wb.Worksheets.Item("WVR")
# Row 20 (idx 29)
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 54 (idx 30)
$ws.Range("H54").Value = 37538
$ws.Range("I54").Value = 34999
$ws.Range("K54").Value = 34999
$ws.Range("M54").Value = -34479

# Row 110 (idx 31)
$ws.Range("H110").Value = 20000
$ws.Range("J110").Value = 20000
$ws.Range("L110").Value = 20000
$ws.Range("N110").Value = -28180
